$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New GUID values (column K, rows 21-137) in top-to-bottom order, matching the
# refactored reference-data ingestion UUIDs from the commit.
$newValues = @(
    "c7c69ba0-f858-4e90-aa64-a6dcd57a8fcb",
    "0e352222-93c8-4557-a641-9158c83c6477",
    "202d2b17-8cf8-4215-a471-c9c07fe2f42e",
    "a606f389-719e-43ff-9785-e23ebb0511a1",
    "22714e1c-6534-4305-bb9f-1d27e47b7e42",
    "51b35066-9c15-4914-bc9c-514d19b09dd5",
    "a896ecfc-80df-41e8-acd5-bc2653941ff5",
    "3b6941b2-e0be-4591-9231-c35bd19e2c92",
    "466b36f1-86ad-4d4c-b856-a1689bd39f77",
    "6e0e1568-9aac-444b-8814-b460028af9e8",
    "6268d063-41cd-4f4e-93eb-7f8680b0a253",
    "9cb81b78-f5cf-4d73-80d3-a09470b3e4c1",
    "bcb64404-e307-41c2-b7ea-a45beb781eb5",
    "97c1251c-860b-42a0-befc-d73e32fb9784",
    "67632208-780a-469c-94e5-00005c442d32",
    "488407d6-5457-4863-a8de-3610e941c486",
    "f8975d05-dcf2-47fc-b4ba-edc92ffe9f09",
    "ab39d8fb-27e8-4d34-866b-d0bfed43c99c",
    "56bdf120-34a3-4546-93ad-53f00a5ea289",
    "0dbe55c0-164d-4bd1-b150-9d7ad8f95666",
    "6005387c-55d9-426d-944e-db531a3e5e2b",
    "954981de-5187-4c3f-8e39-77d9af36d9c8",
    "1adad4cb-b05c-4692-b1ee-4ff26c8ba110",
    "41070fbc-c819-4bd2-b840-f9d5a8cebbd8",
    "98062f52-6acf-4df6-bbb3-4e6bcc9673d8",
    "32fd55b4-3397-438d-a9e8-d8df046e4b9b",
    "b68a60b4-6656-45a5-bc5d-ecf7ac3c9d4e",
    "7b5b33f6-b5a1-40c3-afeb-d56a34c7822f",
    "781c1c32-229e-46f1-a11f-bb8086fe4d62",
    "68da8bd6-6a43-446e-9668-07d4e6fed358",
    "dfe986ef-f14b-4521-abcc-9744e57f97b8",
    "2bbe56cc-d76d-4f6a-8e1f-6564c7e74099",
    "630d39ae-d04b-48e4-a540-489be79757a9",
    "f341e321-a51f-4bb4-80d5-b7e15ad06c0d",
    "1e80f7a7-a503-4d0c-9452-789b2378c40b",
    "aa3ecc24-9883-413e-920b-c57a6712f2f9",
    "0203acab-2899-46dc-b7d7-2464587fa1fc",
    "f21f6279-fec9-496c-be0a-2084ab209892",
    "6b71b40a-baca-46c4-ae49-2017d69fe829",
    "34767309-75b8-4710-8d4e-8bc4d2e6b45c",
    "027ddff8-b1b6-47a2-8c51-c897c9089ce3",
    "68697c64-13db-4735-b1ce-3bbea7c08739",
    "dae0cd7d-6942-474b-995d-0dd20acacdc9",
    "720ae89c-daf0-44c4-8745-ab7bf875ac10",
    "b40c8afc-cc1b-4be8-ad63-59321edecf31",
    "31815b6f-c96e-43c9-8fb6-7fcbaaa63f14",
    "e2d46a1a-4e7a-4e8f-8957-2408eb3a60e0",
    "a3fb126f-285d-46f0-89c1-bafe1df05623",
    "97435109-e58d-4e5e-b476-4a2848ff8a7b",
    "2a9773a4-c60d-4dc3-a4f5-cffdaceb6f83",
    "840f4e8a-148c-424a-9ccd-251896ca737d",
    "9eb8c6eb-14bd-48a4-a8df-b84d3a8562fe",
    "ea561f9b-7b77-416a-b970-2882e408a0c5",
    "452efb45-239a-4c2e-84b2-6a733c59183d",
    "6f5aedb4-f2db-4134-8ed9-e85d8d63a002",
    "38dd49b7-a478-47fb-8aae-f7e38b3adc67",
    "d1794036-118e-4e4b-938c-aebdce18dbd3",
    "5cd83bdd-91f1-4c8e-bd85-c78357d0c588",
    "a45db671-f0ab-43dd-9f0d-0f3521dc6335",
    "87e6e272-b9bc-48f5-a758-e3c1808e9aef",
    "ae3b9b1b-4ca5-472b-a8bc-4e6024d12740",
    "9e8819cb-9d7b-4eb4-835d-04248f9502ef",
    "907bcec3-c618-40e0-9494-cf805e290ee8",
    "7fad2beb-db26-4441-9d2d-c5b1cab67ac7",
    "41e23d15-c190-410f-85f6-9e4c62d4f8aa",
    "c31921e3-fdc8-4dd2-8b2e-d421f580c183",
    "e32650be-23c7-4380-befb-8716e8d23196",
    "61b0dbb4-7332-4162-b830-4ff3b696a8f0",
    "8a05fa10-8a9c-4064-9868-52e526fb1e4a",
    "20deaa51-651b-4923-9fba-392e6cb41bc8",
    "d469a47a-5f43-4816-94e4-319c8e9241ea",
    "6cc8d07a-9051-4615-8b46-8f354fd6d186",
    "8e82ef9d-8778-4a87-9565-312af629bee5",
    "19a0239a-aefe-4e3b-a3ee-ce54461b124c",
    "83f16d06-c095-4392-b5fb-d0455f800897",
    "b91b9c08-954b-4498-b2a3-4c56fd09fce8",
    "89eadf09-6d6d-403d-9667-084794ec1103",
    "3c433083-e284-4cec-bdbd-5095b4ae66ac",
    "02bbfea6-1016-466a-b3df-ef24c42f1d7f",
    "c7b08b01-dacd-424e-8711-7779f02486c0",
    "b3abbd32-b6ee-48a4-baf5-c56e6eaa093f",
    "55d7c9cb-f6ea-40d3-b3f1-3e452ec0217c",
    "55324869-d02e-4a88-b841-b7ad372ee148",
    "418f562d-fe9b-44aa-9b00-636f9da04aa8",
    "8e02ad38-bb3e-403c-baf0-a3c476cdfb3d",
    "1bb7ca21-9b29-4af2-ad84-4a42756e061c",
    "2ffd5191-70cd-45c3-b3ac-054369dc5136",
    "0dd6a98f-deb5-492c-9ae6-177877b5092f",
    "a6b6a0ab-d615-45fa-950b-f31cc2be6c01",
    "64a5dcc4-625f-4c09-a535-79ff35cd3d5a",
    "201b7405-3603-44a0-967b-0483a81cb937",
    "94d1af43-19ba-4cf6-b438-7db42772068a",
    "7067fc04-1f8a-41f2-a41f-fca30d872f50",
    "4639ff42-71a3-4bc3-891e-e9a0c3b122cd",
    "9288295b-a41c-4fdf-b488-7c20b467f5e5",
    "4c8f738f-0c83-4350-80b4-0a9ed20d4799",
    "e2b4adb3-0c0f-4411-9ea1-4e4726acaa6b",
    "38038d0e-91a7-469b-8d15-c59962b507b2",
    "ce5ba600-35b8-4697-86d3-d37108e9cb0a",
    "e3e0cf5d-1075-437d-b12d-8326671dccb4",
    "1a80c7fa-39d4-42eb-bb2f-d6da49e527be",
    "b395917a-423e-44d5-8c78-384d04931861",
    "095d52a4-2eb9-486e-89d8-724353ac08e6",
    "9d056191-a091-4af0-b693-e06d96f85bff",
    "7a8bd44d-615a-4cc3-940e-fd214a499d15",
    "148ab7fb-1eb7-4f8b-a32e-5b614b7b4683",
    "eae3afe6-adec-40b6-9b49-01d00c84971d",
    "b4edc9eb-f96d-48d5-beaa-14bdca734be4",
    "be569bbf-0fda-415b-866b-d9e4b9ccaefd",
    "4b5e09fc-9370-436c-aee7-6a5ff5107628",
    "8153de51-2609-446a-98da-55c7ec1aec91",
    "f7d07424-e9bd-4c42-b905-5a4db6cd6336",
    "70e4d562-c243-4a32-bacf-9dc8f86b8f5b",
    "b2ea2d37-4870-423f-8fc2-fb978f2bf09d",
    "94f9454c-a5f6-4912-a012-acfc819000e7",
    "0571fba4-406d-4357-93a2-54eaca25740b",
    "dbdf1164-ed8c-45f4-9223-20343968fe4e"
)

$startRow = 21
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 11).Value = $newValues[$i]
}
